$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before the old footer/note row (old row 112) so the
# note ("*4/8..." shared-string row) shifts down to row 113, and the new
# row 112 inherits the number-formatting styles used by the data rows
# above it.
$ws.Rows.Item(112).Insert() | Out-Null

# Populate the newly inserted row 112 with the latest day's figures.
$ws.Range("A112").Value = 43967
$ws.Range("B112").Value = 191
$ws.Range("C112").Value = 37481
$ws.Range("D112").Value = 0
$ws.Range("E112").Value = 7584

# Keep the view's active selection on the (now shifted) footer row.
$ws.Range("B113").Select() | Out-Null

# Grow the sheet's print area by one row to keep the footer row included.
$nm = $wb.Names.Item(1)
$nm.RefersTo = "=" + $ws.Name + '!$A$1:$E$114'
